$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2010-2019) to append below the existing table (rows 2-26).
$data = @(
    @(2010, 245.81, 228.77, 17.1, 697.9),
    @(2011, 239.2, 208.74, 17.4, 782.6),
    @(2012, 114.6, 108.07, 17, 568.9),
    @(2013, 105, 100.55, 18.1, 697.8),
    @(2014, 104.2, 93.86, 17.7, 615.6),
    @(2015, 107.4, 92.5, 17.8, 491.8),
    @(2016, 86.9, 72.21, 17.9, 655.9),
    @(2017, 85.4, 84.31, 18, 488),
    @(2018, 74.75, 69.9, 17.8, 786.2),
    @(2019, 76.69, 69.35, 18.8, 565.3)
)

$row = 27
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Formula = "=LN(B$row)"
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Formula = "=LN(D$row)"
    $ws.Cells.Item($row, 6).Value = $r[3]
    $ws.Cells.Item($row, 7).Formula = "=LN(F$row)"
    $ws.Cells.Item($row, 8).Value = $r[4]
    $ws.Cells.Item($row, 9).Formula = "=LN(H$row)"
    $row = $row + 1
}

# Row 1 header height shrinks back to a single-line height once more rows exist.
$ws.Rows.Item(1).RowHeight = 15.75

# Update the selection to the new last cell (also drops the stale topLeftCell scroll anchor).
[void]$ws.Range("I36").Select()
